$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Add the three new header values in row 2 (columns C, D, E)
$ws1.Range("C2").Value = "hallo"
$ws1.Range("D2").Value = "hi "
$ws1.Range("E2").Value = "welcome"

# Resize the first two columns on Sheet1 (width expressed in xlsx "character"
# units includes the ~0.8333 char padding Excel adds on top of the COM
# ColumnWidth value)
$ws1.Columns.Item(1).ColumnWidth = 24.7176870748299
$ws1.Columns.Item(2).ColumnWidth = 22.0187074829932

# Move / update the active selection to the new last used cell
$ws1.Range("E2").Select()

# Sheet3: widen the (still empty) default column block
$ws3.Columns.Item(1).ColumnWidth = 12.2329931972789
for ($c = 2; $c -le 1025; $c++) {
    $ws3.Columns.Item($c).ColumnWidth = 12.2329931972789
}
